$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data originally lived in columns B:C (Semestre / Receita),
# with column A left empty. Deleting column A shifts every cell one
# column to the left, landing the data in A:B - matching the re-uploaded
# workbook where the table now starts at column A instead of B.
$ws.Columns("A").Delete() | Out-Null

# Restore the user's on-screen selection as captured in the re-uploaded
# workbook.
$ws.Range("I27").Select() | Out-Null
